$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.522.45'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.446.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.83%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.445.77'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.87%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.81'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.26%  '
$ws.Range('E11').Value = '  +1.57%  '
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.033.29'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.72%  '
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('E16').Value = '  +0.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.437.34'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.66%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.593.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = '  +8.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '389.24'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.570'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.588.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.83'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.97%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  +0.76%  '
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.182'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.88%  '
$ws.Range('E30').Value = '  +3.52%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('E32').Value = '  -13.37%  '
$ws.Range('E33').Value = '  +1.11%  '
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '24.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.31'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.64%  '
$ws.Range('E38').Value = '  +3.33%  '
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '165.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('E41').Value = '  +4.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.794'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.51'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '42.29'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.616.97'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.11'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.22'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.93%  '
